$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (row 1) ---
$ws.Range("E1").Value = "Job"
$ws.Range("F1").Value = "Location"
$ws.Range("G1").Value = "Phones"
$ws.Range("H1").Value = "Emails"

# --- Last Name (column C) updates ---
$ws.Range("C3").Value = "Higgins MCIOB"
$ws.Range("C4").Value = "Yohanis MCIOB"
$ws.Range("C5").Value = "McLogan CMIOSH LL.M"
$ws.Range("C10").Value = "Gray GradIOSH"
$ws.Range("C11").Value = "Laverty. MCIOB"
$ws.Range("C15").Value = "Gorman (she/her)"
$ws.Range("C19").Value = "Salandy    BSc. (Hons.) GradIOSH"

# --- Email (column H) updates: expand multi-email lists ---
$ws.Range("H6").Value = "lorcan.mulvey@mcaleer-rushe.co.uk , lorcanmulvey@yahoo.ie , lorcan.mulvey@yahoo.ie , lorcan.mulvey@berkeleygroup.co.uk"
$ws.Range("H10").Value = "leergray3@hotmail.co.uk , lee.gray@mcaleer-rushe.co.uk"
$ws.Range("H16").Value = "connor.graham@patton.co.uk , connor.graham@mcaleer-rushe.co.uk"
$ws.Range("H17").Value = "cathal.magee@mcaleer-rushe.co.uk , cathal.magee1@hotmail.co.uk"

# --- Email (column H) updates: clear contents for these rows ---
$ws.Range("H2").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("H12").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("H14").ClearContents()
$ws.Range("H15").ClearContents()
$ws.Range("H19").ClearContents()
$ws.Range("H20").ClearContents()
$ws.Range("H21").ClearContents()
$ws.Range("H22").ClearContents()
